$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the two rows whose records were dropped -------------
# Row 26 ("RM 232") is removed; everything below shifts up one row.
$ws.Range("A26:F26").EntireRow.Delete()
# After the shift above, the old "SC 92" record (previously row 28) is now
# row 27; remove it too so the remaining SC rows shift up again.
$ws.Range("A27:F27").EntireRow.Delete()

# --- Step 2: per-cell value corrections on the (now renumbered) sheet ---
$ws.Range("F2").Value = ""

$ws.Range("F5").Value = 17.66

$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43

$ws.Range("E8").Value = ""

$ws.Range("F10").Value = ""

$ws.Range("E12").Value = -5.3

$ws.Range("F13").Value = ""

$ws.Range("E14").Value = ""

$ws.Range("E17").Value = -7.3

$ws.Range("E18").Value = -8.5

$ws.Range("E19").Value = ""

$ws.Range("E20").Value = ""

$ws.Range("E23").Value = -7

$ws.Range("F24").Value = 16.78

# row 27 is now "SC 101" (was old row 29)
$ws.Range("B27").Value = -20.4
$ws.Range("E27").Value = ""

# row 28 is now "SC 105" (was old row 30)
$ws.Range("F28").Value = ""

# row 29 is now "SC 119" (was old row 31)
$ws.Range("B29").Value = ""

# row 30 is now "SC 120" (was old row 32)
$ws.Range("F30").Value = 16.89

# row 32 is now "SC 193" (was old row 34)
$ws.Range("B32").Value = ""
